# Applies the recalculated VESIcal output values (Temp column + all downstream
# dissolved/eqfluid/eqfluid_wtemps/SaturationPs/SatPs_wtemps columns) to match
# the updated model run referenced in the commit.
$wb = $excel.ActiveWorkbook

# --- Original_User_Data ---
$ws = $wb.Worksheets.Item("Original_User_Data")
$ws.Range("U2").Value = 1299.09471
$ws.Range("U3").Value = 1283.41999
$ws.Range("U4").Value = 1255.15376

# --- dissolved ---
$ws = $wb.Worksheets.Item("dissolved")
$ws.Range("U2").Value = 1299.09471
$ws.Range("U3").Value = 1283.41999
$ws.Range("U4").Value = 1255.15376

# --- eqfluid ---
$ws = $wb.Worksheets.Item("eqfluid")
$ws.Range("U2").Value = 1299.09471
$ws.Range("U3").Value = 1283.41999
$ws.Range("U4").Value = 1255.15376

# --- eqfluid_wtemps ---
$ws = $wb.Worksheets.Item("eqfluid_wtemps")
$ws.Range("U2").Value = 1299.09471
$ws.Range("V2").Value = 0.474060274275901
$ws.Range("W2").Value = 0.525939725724099
$ws.Range("U3").Value = 1283.41999
$ws.Range("U4").Value = 1255.15376

# --- SaturationPs ---
$ws = $wb.Worksheets.Item("SaturationPs")
$ws.Range("U2").Value = 1299.09471
$ws.Range("V2").Value = 60
$ws.Range("X2").Value = 0.469913142701622
$ws.Range("Y2").Value = 0.530086857298378
$ws.Range("Z2").Value = 0.000835660852316802
$ws.Range("AA2").Value = 0.0008356608523168019
$ws.Range("U3").Value = 1283.41999
$ws.Range("Z3").Value = 0.0000375646222415579
$ws.Range("AA3").Value = 0.00003756462224155789
$ws.Range("U4").Value = 1255.15376
$ws.Range("V4").Value = 100
$ws.Range("X4").Value = 0.29235352352064
$ws.Range("Y4").Value = 0.70764647647936
$ws.Range("Z4").Value = 0.000634775486311526
$ws.Range("AA4").Value = 0.0006347754863115258
$ws.Range("V5").Value = 2500
$ws.Range("X5").Value = 0.796513677120962
$ws.Range("Y5").Value = 0.203486322879038
$ws.Range("Z5").Value = 0.00123226570304335
$ws.Range("AA5").Value = 0.00123226570304335
$ws.Range("X6").Value = 0.836894576666795
$ws.Range("Y6").Value = 0.163105423333205
$ws.Range("Z6").Value = 0.000226271338451074
$ws.Range("AA6").Value = 0.0002262713384510739
$ws.Range("Z7").Value = 0.0129029986829266
$ws.Range("AA7").Value = 0.01290299868292659
$ws.Range("Z8").Value = 0.00105217920776629
$ws.Range("AA8").Value = 0.00105217920776629
$ws.Range("V9").Value = 2540
$ws.Range("Z9").Value = 0.0160925216776623
$ws.Range("AA9").Value = 0.01609252167766231
$ws.Range("V10").Value = 1100
$ws.Range("X10").Value = 0.972472424256134
$ws.Range("Y10").Value = 0.0275275757438659
$ws.Range("Z10").Value = 0.007924265844899119
$ws.Range("AA10").Value = 0.007924265844899111
$ws.Range("V11").Value = 1790
$ws.Range("X11").Value = 0.972874630367114
$ws.Range("Y11").Value = 0.0271253696328862
$ws.Range("Z11").Value = 0.00667100825932183
$ws.Range("AA11").Value = 0.006671008259321829
$ws.Range("V12").Value = 1730
$ws.Range("X12").Value = 0.97561394228405
$ws.Range("Y12").Value = 0.0243860577159497
$ws.Range("Z12").Value = 0.008636970211907451
$ws.Range("AA12").Value = 0.008636970211907449
$ws.Range("V13").Value = 2090
$ws.Range("X13").Value = 0.951891002202502
$ws.Range("Y13").Value = 0.048108997797498
$ws.Range("Z13").Value = 0.00294119534061353
$ws.Range("AA13").Value = 0.002941195340613532
$ws.Range("V14").Value = 1730
$ws.Range("X14").Value = 0.950741498987455
$ws.Range("Y14").Value = 0.0492585010125445
$ws.Range("Z14").Value = 0.00286390060466055
$ws.Range("AA14").Value = 0.002863900604660548
$ws.Range("X15").Value = 0.231707821016213
$ws.Range("Y15").Value = 0.768292178983787
$ws.Range("Z15").Value = 0.00009309173835698601
$ws.Range("AA15").Value = 0.00009309173835698606
$ws.Range("X16").Value = 0.456749735218006
$ws.Range("Y16").Value = 0.543250264781994
$ws.Range("Z16").Value = 0.000938033056325138
$ws.Range("AA16").Value = 0.0009380330563251388
$ws.Range("X17").Value = 0.684728944146793
$ws.Range("Y17").Value = 0.315271055853207
$ws.Range("Z17").Value = 0.000430867551518525
$ws.Range("AA17").Value = 0.0004308675515185249

# --- SatPs_wtemps ---
$ws = $wb.Worksheets.Item("SatPs_wtemps")
$ws.Range("U2").Value = 1299.09471
$ws.Range("V2").Value = 60
$ws.Range("W2").Value = 0.493184396428161
$ws.Range("X2").Value = 0.506815603571839
$ws.Range("Y2").Value = 0.000609818056595781
$ws.Range("Z2").Value = 0.000609818056595781
$ws.Range("U3").Value = 1283.41999
$ws.Range("V3").Value = 110
$ws.Range("W3").Value = 0.266594556578593
$ws.Range("X3").Value = 0.733405443421407
$ws.Range("Y3").Value = 0.000699690455656429
$ws.Range("Z3").Value = 0.0006996904556564294
$ws.Range("U4").Value = 1255.15376
$ws.Range("V4").Value = 90
$ws.Range("W4").Value = 0.337738175495508
$ws.Range("X4").Value = 0.662261824504492
$ws.Range("Y4").Value = 0.000807422690952036
$ws.Range("Z4").Value = 0.0008074226909520355
$ws.Range("V5").Value = 2540
$ws.Range("W5").Value = 0.817547776604673
$ws.Range("X5").Value = 0.182452223395327
$ws.Range("Y5").Value = 0.0015319426079297
$ws.Range("Z5").Value = 0.0015319426079297
$ws.Range("W6").Value = 0.8552141507186199
$ws.Range("X6").Value = 0.14478584928138
$ws.Range("Y6").Value = 0.000849378234060265
$ws.Range("Z6").Value = 0.0008493782340602644
$ws.Range("Y7").Value = 0.00344178052689819
$ws.Range("Z7").Value = 0.003441780526898187
$ws.Range("V8").Value = 1650
$ws.Range("Y8").Value = 0.0152803231983435
$ws.Range("Z8").Value = 0.01528032319834351
$ws.Range("Y9").Value = 0.008153007612025051
$ws.Range("Z9").Value = 0.008153007612025049
$ws.Range("V10").Value = 1090
$ws.Range("W10").Value = 0.972915709108933
$ws.Range("X10").Value = 0.0270842908910666
$ws.Range("Y10").Value = 0.00885469016241785
$ws.Range("Z10").Value = 0.008854690162417838
$ws.Range("V11").Value = 1780
$ws.Range("W11").Value = 0.973132721034933
$ws.Range("X11").Value = 0.0268672789650666
$ws.Range("Y11").Value = 0.00591552141542247
$ws.Range("Z11").Value = 0.005915521415422469
$ws.Range("V12").Value = 1720
$ws.Range("W12").Value = 0.975859890544081
$ws.Range("X12").Value = 0.0241401094559188
$ws.Range("Y12").Value = 0.00808765212115037
$ws.Range("Z12").Value = 0.008087652121150367
$ws.Range("V14").Value = 1730
$ws.Range("W14").Value = 0.951016786572314
$ws.Range("X14").Value = 0.048983213427686
$ws.Range("Y14").Value = 0.00335004855562078
$ws.Range("Z14").Value = 0.003350048555620779
$ws.Range("V15").Value = 1280
$ws.Range("W15").Value = 0.228644415713756
$ws.Range("X15").Value = 0.771355584286244
$ws.Range("Y15").Value = 0.00147492657477965
$ws.Range("Z15").Value = 0.00147492657477965
$ws.Range("V16").Value = 4910
$ws.Range("W16").Value = 0.458904496315431
$ws.Range("X16").Value = 0.541095503684569
$ws.Range("Y16").Value = 0.00176716229182296
$ws.Range("Z16").Value = 0.001767162291822961
$ws.Range("V17").Value = 1590
$ws.Range("W17").Value = 0.679642799711047
$ws.Range("X17").Value = 0.320357200288953
$ws.Range("Y17").Value = 0.00191408285739597
$ws.Range("Z17").Value = 0.00191408285739597
